$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: "Championed GenAI-powered auto-extraction capabilities..." bullet.
# Merge the "GenAI" run and the "-powered auto-extraction capabilities" run
# (currently split apart by <w:proofErr w:type="spellStart"/> /
# <w:proofErr w:type="spellEnd"/> markers) into a single bold ("Strong")
# run reading "GenAI-powered auto-extraction capabilities", and drop the
# proofErr markers entirely.
# ---------------------------------------------------------------------------

$markerChar = [char]0xE000

$rngGenAI = $d.Content
$rngGenAI.Find.Execute("GenAI", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$genaiStart = $rngGenAI.Start

$rngCap = $d.Content
$rngCap.Start = $genaiStart
$rngCap.End = $d.Content.End
$rngCap.Find.Execute("auto-extraction capabilities", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$capEnd = $rngCap.End

# Extend the start one character to the left (onto the trailing non-breaking
# space of the preceding "Championed " run) so the edit crosses the
# w:proofErr boundary and Word folds it away; a private-use marker
# character is appended so the run is guaranteed to actually change (a
# same-text rewrite is a no-op) and so we can re-select it afterwards.
$fullGenAI = $d.Range($genaiStart - 1, $capEnd)
$fullGenAI.Text = [char]0x00A0 + "GenAI-powered auto-extraction capabilities" + $markerChar

# Re-apply the "Strong" character style to the merged run (the crossing
# edit above adopts the plain "Championed " run's formatting).
$rngStyle = $d.Content
$rngStyle.Find.Execute("GenAI-powered auto-extraction capabilities" + $markerChar, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngStyle.Style = "Strong"

# Strip the temporary marker character back out.
$rngMarker = $d.Content
$rngMarker.Find.Execute($markerChar, $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# ---------------------------------------------------------------------------
# Edit 2: add a new "Professional certificate  program in GenAI and ML from
# IIT Madras (Oct 2025-Apr 2026)" bullet right before the existing
# "Chief Technology Officer, ISB – 2024" bullet under PROFESSIONAL
# CERTIFICATIONS.
# ---------------------------------------------------------------------------

$rngCTO = $d.Content
$rngCTO.Find.Execute("Chief Technology Officer, ISB", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$ctoPara = $rngCTO.Paragraphs(1)
$ctoStart = $ctoPara.Range.Start

$insertPoint = $d.Range($ctoStart, $ctoStart)
$insertPoint.InsertBefore("Professional certificate  program in GenAI and ML from IIT Madras (Oct 2025-Apr 2026)`r")

# ---------------------------------------------------------------------------
# Edit 3: "Certified SAFe Agilist – 2022" bullet. Merge the "Certified ",
# "SAFe" and " Agilist – 2022" runs (split apart by proofErr markers around
# "SAFe") into a single run, dropping the proofErr markers.
# ---------------------------------------------------------------------------

$rngSafe = $d.Content
$rngSafe.Find.Execute("SAFe", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$safePara = $rngSafe.Paragraphs(1)
$safeParaStart = $safePara.Range.Start

$rngAfterSafe = $d.Content
$rngAfterSafe.Start = $safeParaStart
$rngAfterSafe.End = $d.Content.End
$rngAfterSafe.Find.Execute("SAFe Agilist", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$afterSafeEnd = $rngAfterSafe.End

$rngYear = $d.Content
$rngYear.Start = $afterSafeEnd
$rngYear.End = $d.Content.End
$rngYear.Find.Execute("2022", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$yearEnd = $rngYear.End

$fullSafe = $d.Range($safeParaStart, $yearEnd)
$fullSafe.Text = "Certified SAFe Agilist" + [char]0x00A0 + [char]0x2013 + " 2022" + $markerChar

$rngMarker2 = $d.Content
$rngMarker2.Find.Execute($markerChar, $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
